$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 346.66666
$ws.Range("I33").Value = 270.25
$ws.Range("K33").Value = 270.25
$ws.Range("M33").Value = -41.25

$ws.Range("H62").Value = 1985.3
$ws.Range("I62").Value = 1150
$ws.Range("J62").Value = 3934.3333
$ws.Range("K62").Value = 1150
$ws.Range("L62").Value = 3934.3333
$ws.Range("M62").Value = -526
$ws.Range("N62").Value = -5182.3333

$ws.Range("H65").Value = 1985.3
$ws.Range("I65").Value = 1150
$ws.Range("J65").Value = 3934.3333
$ws.Range("K65").Value = 5750
$ws.Range("L65").Value = 19671.6665
$ws.Range("M65").Value = -2630
$ws.Range("N65").Value = -25911.6665

$ws.Range("H88").Value = 5685.5713
$ws.Range("I88").Value = 2384.8572
$ws.Range("K88").Value = 2384.8572
$ws.Range("M88").Value = -1978.8572

$ws.Range("H91").Value = 5685.5713
$ws.Range("I91").Value = 2384.8572
$ws.Range("K91").Value = 2384.8572
$ws.Range("M91").Value = -980.8571999999999

$ws.Range("H100").Value = 33334800
$ws.Range("I100").Value = 40000760
$ws.Range("K100").Value = 40000760
$ws.Range("M100").Value = -40000219

$ws.Range("H106").Value = 3261.6
$ws.Range("I106").Value = 2290.6667
$ws.Range("K106").Value = 2290.6667
$ws.Range("M106").Value = -1659.6667

$ws.Range("H116").Value = 534515.9
$ws.Range("I116").Value = 2001998
$ws.Range("J116").Value = 10415.143
$ws.Range("K116").Value = 2001998
$ws.Range("L116").Value = 10415.143
$ws.Range("M116").Value = -1998556
$ws.Range("N116").Value = -17299.143

$ws.Range("H129").Value = 913.24
$ws.Range("J129").Value = 966.3111
$ws.Range("L129").Value = 2898.9333
$ws.Range("N129").Value = -12898.9333

$ws.Range("H138").Value = 3098.798
$ws.Range("J138").Value = 4271.4463
$ws.Range("L138").Value = 12814.3389
$ws.Range("N138").Value = -23094.3389

$ws.Range("H141").Value = 6642.564
$ws.Range("J141").Value = 3064.1667
$ws.Range("L141").Value = 9192.500100000001
$ws.Range("N141").Value = -19552.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2999
$ws.Range("I2").Value = 2999
$ws.Range("K2").Value = 2999
$ws.Range("M2").Value = -2886

$ws.Range("H32").Value = 4879.3784
$ws.Range("I32").Value = 3584
$ws.Range("J32").Value = 9222.706
$ws.Range("K32").Value = 3584
$ws.Range("L32").Value = 9222.706
$ws.Range("M32").Value = -3297
$ws.Range("N32").Value = -9796.706

$ws.Range("H116").Value = 2999
$ws.Range("I116").Value = 2999
$ws.Range("K116").Value = 2999
$ws.Range("M116").Value = -705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2999
$ws.Range("I3").Value = 2999
$ws.Range("K3").Value = 2999
$ws.Range("M3").Value = -2885

$ws.Range("H99").Value = 3139.2
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 3549
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 3549
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -6545

$ws.Range("H105").Value = 1676.9722
$ws.Range("I105").Value = 1679.1
$ws.Range("J105").Value = 1666.3334
$ws.Range("K105").Value = 1679.1
$ws.Range("L105").Value = 1666.3334
$ws.Range("M105").Value = 67.90000000000009
$ws.Range("N105").Value = -5160.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10640692
$ws.Range("I31").Value = 1323.5333
$ws.Range("J31").Value = 29416050
$ws.Range("K31").Value = 1323.5333
$ws.Range("L31").Value = 29416050
$ws.Range("M31").Value = -1028.5333
$ws.Range("N31").Value = -29416640

$ws.Range("H34").Value = 10640692
$ws.Range("I34").Value = 1323.5333
$ws.Range("J34").Value = 29416050
$ws.Range("K34").Value = 1323.5333
$ws.Range("L34").Value = 29416050
$ws.Range("M34").Value = -1121.5333
$ws.Range("N34").Value = -29416454

$ws.Range("H99").Value = 11770584
$ws.Range("I99").Value = 20004242
$ws.Range("J99").Value = 8214.286
$ws.Range("K99").Value = 20004242
$ws.Range("L99").Value = 8214.286
$ws.Range("M99").Value = -20002744
$ws.Range("N99").Value = -11210.286

$ws.Range("H105").Value = 2051.5557
$ws.Range("I105").Value = 1880.7693
$ws.Range("K105").Value = 1880.7693
$ws.Range("M105").Value = -133.7692999999999

$ws.Range("H126").Value = 11770584
$ws.Range("I126").Value = 20004242
$ws.Range("J126").Value = 8214.286
$ws.Range("K126").Value = 60012726
$ws.Range("L126").Value = 24642.858
$ws.Range("M126").Value = -60010256
$ws.Range("N126").Value = -29582.858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7937271.5
$ws.Range("J131").Value = 863.7059
$ws.Range("L131").Value = 2591.1177
$ws.Range("N131").Value = -12671.1177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6447.0713
$ws.Range("I70").Value = 5895.6665
$ws.Range("J70").Value = 8101.2856
$ws.Range("K70").Value = 5895.6665
$ws.Range("L70").Value = 8101.2856
$ws.Range("M70").Value = -5625.6665
$ws.Range("N70").Value = -8641.285599999999

$ws.Range("H73").Value = 6447.0713
$ws.Range("I73").Value = 5895.6665
$ws.Range("J73").Value = 8101.2856
$ws.Range("K73").Value = 5895.6665
$ws.Range("L73").Value = 8101.2856
$ws.Range("M73").Value = -4959.6665
$ws.Range("N73").Value = -9973.285599999999

$ws.Range("H97").Value = 730
$ws.Range("I97").Value = 730
$ws.Range("K97").Value = 730
$ws.Range("M97").Value = -234

$ws.Range("H113").Value = 38104.332
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 56156.5
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 56156.5
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -60496.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 19999.5
$ws.Range("J29").Value = 19999.5
$ws.Range("L29").Value = 19999.5
$ws.Range("N29").Value = -20589.5

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1298
$ws.Range("N61").Value = -2904

$ws.Range("H93").Value = 6946429.5
$ws.Range("I93").Value = 12347167
$ws.Range("J93").Value = 2624.1428
$ws.Range("K93").Value = 12347167
$ws.Range("L93").Value = 2624.1428
$ws.Range("M93").Value = -12345919
$ws.Range("N93").Value = -5120.1428

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H113").Value = 443.86957
$ws.Range("I113").Value = 342.35294
$ws.Range("J113").Value = 731.5
$ws.Range("K113").Value = 1027.05882
$ws.Range("L113").Value = 2194.5
$ws.Range("M113").Value = 1142.94118
$ws.Range("N113").Value = -6534.5

Write-Host "Updated pricing/profit values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
